# Updated cryptos list on Fri Nov  1 13:32:18 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.259.00"
$ws.Range("E2").Value = "  -2.39%  "
$ws.Range("D3").Value = "2.536.67"
$ws.Range("E3").Value = "  -3.29%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "578.45"
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("D6").Value = "170.08"
$ws.Range("E6").Value = "  -2.05%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "0.518"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "2.535.66"
$ws.Range("E9").Value = "  -3.25%  "
$ws.Range("E10").Value = "  -4.29%  "
$ws.Range("E11").Value = "  -1.29%  "
$ws.Range("E12").Value = "  -1.57%  "
$ws.Range("D13").Value = "4.93"
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("D14").Value = "2.987.39"
$ws.Range("E14").Value = "  -3.88%  "
$ws.Range("D15").Value = "70.140.46"
$ws.Range("E15").Value = "  -2.40%  "
$ws.Range("E16").Value = "  -6.25%  "
$ws.Range("D17").Value = "25.32"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("D18").Value = "2.533.76"
$ws.Range("E18").Value = "  -3.77%  "
$ws.Range("D19").Value = "'7.90"
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("E20").Value = "  -6.18%  "
$ws.Range("D21").Value = "353.14"
$ws.Range("E21").Value = "  -5.40%  "
$ws.Range("E22").Value = "  -2.80%  "
$ws.Range("D23").Value = "2.02"
$ws.Range("E23").Value = "  +0.88%  "
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").Value = "70.04"
$ws.Range("E25").Value = "  -1.62%  "
$ws.Range("D26").Value = "4.03"
$ws.Range("E26").Value = "  -3.47%  "
$ws.Range("E27").Value = "  -1.71%  "
$ws.Range("D28").Value = "2.657.86"
$ws.Range("E29").Value = "  +0.75%  "
$ws.Range("D30").Value = "0.0₃0917"
$ws.Range("E30").Value = "  -2.63%  "
$ws.Range("D31").Value = "7.92"
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "469.94"
$ws.Range("E32").Value = "  -2.72%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "1.27"
$ws.Range("E33").Value = "  -3.29%  "
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  +2.97%  "
$ws.Range("D37").Value = "157.55"
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("D38").Value = "19.03"
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("D39").Value = "'18.60"
$ws.Range("E39").Value = "  -3.09%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("E42").Value = "  -0.96%  "
$ws.Range("D43").Value = "1.61"
$ws.Range("E43").Value = "  -6.08%  "
$ws.Range("E44").Value = "  -8.06%  "
$ws.Range("E45").Value = "  -13.79%  "
$ws.Range("D46").Value = "38.33"
$ws.Range("E46").Value = "  -1.82%  "
$ws.Range("D47").Value = "144.66"
$ws.Range("E47").Value = "  -3.24%  "
$ws.Range("D48").Value = "0.536"
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("E49").Value = "  -2.91%  "
$ws.Range("E50").Value = "  -2.71%  "
$ws.Range("E51").Value = "  -0.58%  "
